# Tester for final approach
# - Populates sheet "d2" with benchmark data (matching the layout/style of "d1")
# - Moves the active tab/selection from "VM" to "d2"
# - Updates the selection remembered on "d1"

$wb = $excel.ActiveWorkbook

$wsVM = $wb.Worksheets.Item("VM")
$wsD1 = $wb.Worksheets.Item("d1")
$wsD2 = $wb.Worksheets.Item("d2")

# --- d2: bring in the same number formatting/font used by the other data sheets ---
# (copy format only, so the shared style table / fonts stay untouched)
$wsD1.Range("A1:C25").Copy()
$wsD2.Range("A1:C25").PasteSpecial(-4122)   # xlPasteFormats

# --- d2: fill in the benchmark values ---
$data = @(
    @(0,     0,     0),
    @(0,     0,     0),
    @(1909,  6230,  0),
    @(0,     0,     0),
    @(0,     0,     0),
    @(0,     0,     0),
    @(0,     0,     0),
    @(0,     0,     0),
    @(0,     0,     0),
    @(4039,  23375, 0),
    @(0,     0,     0),
    @(0,     0,     0),
    @(0,     0,     0),
    @(0,     0,     0),
    @(0,     0,     0),
    @(0,     0,     0),
    @(0,     0,     0),
    @(0,     0,     0),
    @(0,     0,     0),
    @(11384, 0,     0),
    @(0,     0,     0),
    @(0,     0,     0),
    @(0,     0,     0),
    @(0,     0,     0),
    @(0,     0,     0)
)

for ($r = 1; $r -le 25; $r++) {
    for ($c = 1; $c -le 3; $c++) {
        $wsD2.Cells.Item($r, $c).Value = $data[$r - 1][$c - 1]
    }
}

# --- d1: remembered selection moves from D30 to D34 ---
$wsD1.Range("D34").Select()

# --- d2 becomes the active sheet/tab, with A10 selected ---
$wsD2.Activate()
$wsD2.Range("A10").Select()
